$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '63.533.18'
$ws.Range("E2").Value = '  +0.35%  '

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '3.099.50'
$ws.Range("E3").Value = '  -0.70%  '

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '''582.95'
$ws.Range("E5").Value = '  -0.56%  '

$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D6").Value = '''145.03'
$ws.Range("E6").Value = '  -0.04%  '

$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").Value = '3.092.99'
$ws.Range("E8").Value = '  -0.66%  '

$ws.Range("B9").Value = 'XRP'
$ws.Range("C9").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D9").Value = '''0.528'
$ws.Range("E9").Value = '  -0.33%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '''0.161'
$ws.Range("E10").Value = '  +6.85%  '

$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").Value = '''5.62'
$ws.Range("E11").Value = '  -2.74%  '

$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").Value = '''0.456'
$ws.Range("E12").Value = '  -2.51%  '

$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D13").Value = '''0.0000245'
$ws.Range("E13").Value = '  -1.31%  '

$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = '''37.10'
$ws.Range("E14").Value = '  +4.01%  '

$ws.Range("B15").Value = 'TRON'
$ws.Range("C15").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D15").Value = '''0.121'
$ws.Range("E15").Value = '  -1.22%  '

$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").Value = '3.609.85'
$ws.Range("E16").Value = '  -0.69%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '63.342.15'
$ws.Range("E17").Value = '  +0.21%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.096.39'
$ws.Range("E18").Value = '  -0.69%  '

$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").Value = '''7.06'
$ws.Range("E19").Value = '  -1.61%  '

$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '''460.29'
$ws.Range("E20").Value = '  -1.62%  '

$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").Value = '''14.21'
$ws.Range("E21").Value = '  +0.94%  '

$ws.Range("B22").Value = 'Polygon'
$ws.Range("C22").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D22").Value = '''0.722'
$ws.Range("E22").Value = '  -0.59%  '

$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").Value = '''7.41'
$ws.Range("E23").Value = '  -1.75%  '

$ws.Range("B24").Value = 'InternetComputer(DFINITY)'
$ws.Range("C24").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D24").Value = '''12.90'
$ws.Range("E24").Value = '  -3.04%  '

$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = '''81.07'
$ws.Range("E25").Value = '  -1.26%  '

$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").Value = '''2.19'
$ws.Range("E26").Value = '  +0.49%  '

$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").Value = '''1.00'
$ws.Range("E27").Value = '  -0.05%  '

$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Value = '''9.14'
$ws.Range("E28").Value = '  +9.40%  '

$ws.Range("B29").Value = 'FirstDigitalUSD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  -0.01%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = '''2.66'
$ws.Range("E30").Value = '  -0.66%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = '''2.19'
$ws.Range("E31").Value = '  -2.23%  '

$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").Value = '''6.92'
$ws.Range("E32").Value = '  +0.82%  '

$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '''26.63'
$ws.Range("E33").Value = '  -1.36%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '''0.109'
$ws.Range("E34").Value = '  -0.59%  '

$ws.Range("B35").Value = 'PEPE'
$ws.Range("C35").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D35").Value = '0.0₃0847'
$ws.Range("E35").Value = '  -2.81%  '

$ws.Range("B36").Value = 'dogwifhat'
$ws.Range("C36").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D36").Value = '''3.39'
$ws.Range("E36").Value = '  +2.31%  '

$ws.Range("B37").Value = 'Mantle'
$ws.Range("C37").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D37").Value = '''1.03'
$ws.Range("E37").Value = '  -0.77%  '

$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").Value = '''2.31'
$ws.Range("E38").Value = '  -5.22%  '

$ws.Range("B39").Value = 'Filecoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D39").Value = '''5.98'
$ws.Range("E39").Value = '  -0.70%  '

$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").Value = '''50.25'
$ws.Range("E40").Value = '  -1.30%  '

$ws.Range("B41").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D41").Value = '''436.09'
$ws.Range("E41").Value = '  +0.62%  '

$ws.Range("B42").Value = 'Cosmos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D42").Value = '''8.67'
$ws.Range("E42").Value = '  -0.67%  '

$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '2.883.38'
$ws.Range("E43").Value = '  -1.61%  '

$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '''0.0366'
$ws.Range("E44").Value = '  -1.01%  '

$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").Value = '''0.271'
$ws.Range("E45").Value = '  -3.15%  '

$ws.Range("B46").Value = 'Kaspa'
$ws.Range("C46").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D46").Value = '''0.108'
$ws.Range("E46").Value = '  -3.31%  '

$ws.Range("B47").Value = 'Arweave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D47").Value = '''36.23'
$ws.Range("E47").Value = '  +2.55%  '

$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").Value = '''124.58'
$ws.Range("E48").Value = '  +1.08%  '

$ws.Range("B49").Value = 'USDe'
$ws.Range("C49").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D49").Value = '''0.999'
$ws.Range("E49").Value = '  -0.03%  '

$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").Value = '''0.109'
$ws.Range("E50").Value = '  -1.49%  '

$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").Value = '''23.99'
$ws.Range("E51").Value = '  -2.47%  '
